$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.379.62"
$ws.Range("E2").Value = "  -0.74%  "

# Row 3
$ws.Range("D3").Value = "1.708.64"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9951"
$ws.Range("E4").Value = "  -0.47%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.12"
$ws.Range("E5").Value = "  -3.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9971"
$ws.Range("E6").Value = "  -0.27%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4883"
$ws.Range("E7").Value = "  -0.95%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2581"
$ws.Range("E8").Value = "  -3.88%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06170"
$ws.Range("E9").Value = "  -1.96%  "

# Row 10
$ws.Range("D10").Value = "1.715.32"
$ws.Range("E10").Value = "  -1.39%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06949"
$ws.Range("E11").Value = "  -1.50%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.51"
$ws.Range("E12").Value = "  -1.54%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5979"
$ws.Range("E13").Value = "  -2.82%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.453"
$ws.Range("E14").Value = "  -2.98%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "76.44"
$ws.Range("E15").Value = "  -1.97%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9968"
$ws.Range("E16").Value = "  -0.28%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "26.241.77"
$ws.Range("E17").Value = "  -1.29%  "

# Row 18
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9949"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007068"
$ws.Range("E19").Value = "  -4.04%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("E20").Value = "  -3.12%  "

# Row 21
$ws.Range("D21").Value = "1.936.68"
$ws.Range("E21").Value = "  -1.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.377"
$ws.Range("E22").Value = "  -4.93%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.376"
$ws.Range("E23").Value = "  -4.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.012"
$ws.Range("E24").Value = "  -4.73%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.23"
$ws.Range("E25").Value = "  -2.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.15"
$ws.Range("E26").Value = "  -2.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.404"
$ws.Range("E27").Value = "  -1.16%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.726"
$ws.Range("E28").Value = "  -2.22%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "105.32"
$ws.Range("E29").Value = "  -2.38%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.876"
$ws.Range("E30").Value = "  -4.34%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07932"
$ws.Range("E31").Value = "  -1.65%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.591"
$ws.Range("E32").Value = "  -3.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04445"
$ws.Range("E33").Value = "  -3.86%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.605"
$ws.Range("E34").Value = "  -0.24%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9906"
$ws.Range("E35").Value = "  -2.82%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6143"
$ws.Range("E36").Value = "  -3.75%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9478"
$ws.Range("E37").Value = "  +5.33%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.991"
$ws.Range("E38").Value = "  -2.44%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.365"
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9971"
$ws.Range("E40").Value = "  -0.76%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01476"
$ws.Range("E41").Value = "  -1.96%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.45"
$ws.Range("E42").Value = "  -2.44%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.399"
$ws.Range("E43").Value = "  -0.27%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3796"
$ws.Range("E44").Value = "  -3.50%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.804"
$ws.Range("E45").Value = "  -1.55%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1146"
$ws.Range("E46").Value = "  -3.53%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05341"
$ws.Range("E47").Value = "  -1.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.39"
$ws.Range("E48").Value = "  -0.73%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.649"
$ws.Range("E49").Value = "  -2.64%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "51.02"
$ws.Range("E50").Value = "  -1.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9993"
$ws.Range("E51").Value = "  -0.33%  "
